# Correction de doublons dans l'Aude
# The "Aude" department row (row 3, dpt=11) had its weekly flux-vision
# figures (sem_01..sem_53, columns I..BI) double-counted because of a
# duplicate merge upstream. Halve every weekly value on that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 9; $col -le 61; $col++) {
    $cell = $ws.Cells.Item(3, $col)
    $old = $cell.Value2
    $cell.Value = $old / 2
}

# Widen column B (department name) so "Haute-Garonne" etc. aren't clipped.
$ws.Columns.Item(2).ColumnWidth = 19.8

# Leave the cursor parked on BB20, matching where the author ended up.
$ws.Range("BB20").Select() | Out-Null
